$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1) - new columns I and J, styled like the other headers (copy format from H1)
$ws.Range("H1").Copy($ws.Range("I1"))
$ws.Range("I1").Value = "I0"
$ws.Range("H1").Copy($ws.Range("J1"))
$ws.Range("J1").Value = "IF"

# Data rows 2-13 for columns I and J
$data = @{
    2  = @(9, 9)
    3  = @(3, 4)
    4  = @(6, 7)
    5  = @(6, 6)
    6  = @(4, 4)
    7  = @(6, 7)
    8  = @(8, 9)
    9  = @(9, 9)
    10 = @(9, 9)
    11 = @(9, 9)
    12 = @(7, 7)
    13 = @(5, 5)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 9).Value = $vals[0]
    $ws.Cells.Item($row, 10).Value = $vals[1]
}
